# Auto-generated Excel COM-interop script
# Updates cached market-price / profit columns (H-N) in the Goblin Profits workbook
# as produced by the scheduled price-update runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 315.5
$ws.Cells.Item(2, 9).Value = 315.5
$ws.Cells.Item(2, 11).Value = 315.5
$ws.Cells.Item(2, 13).Value = -202.5

$ws.Cells.Item(9, 8).Value = 95.333336
$ws.Cells.Item(9, 9).Value = 86.28570999999999
$ws.Cells.Item(9, 11).Value = 86.28570999999999
$ws.Cells.Item(9, 13).Value = 82.71429000000001

$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 13).ClearContents()

$ws.Cells.Item(39, 8).Value = 208.57143
$ws.Cells.Item(39, 9).Value = 33.3
$ws.Cells.Item(39, 11).Value = 99.89999999999999
$ws.Cells.Item(39, 13).Value = 196.1

$ws.Cells.Item(58, 8).Value = 19235694
$ws.Cells.Item(58, 10).Value = 14749.75
$ws.Cells.Item(58, 12).Value = 44249.25
$ws.Cells.Item(58, 14).Value = -44549.25

$ws.Cells.Item(62, 8).Value = 66908
$ws.Cells.Item(62, 9).Value = 91311.55499999999
$ws.Cells.Item(62, 11).Value = 91311.55499999999
$ws.Cells.Item(62, 13).Value = -90687.55499999999

$ws.Cells.Item(65, 8).Value = 66908
$ws.Cells.Item(65, 9).Value = 91311.55499999999
$ws.Cells.Item(65, 11).Value = 456557.775
$ws.Cells.Item(65, 13).Value = -453437.775

$ws.Cells.Item(74, 8).Value = 18430.285
$ws.Cells.Item(74, 9).Value = 18430.285
$ws.Cells.Item(74, 11).Value = 18430.285
$ws.Cells.Item(74, 13).Value = -17494.285

$ws.Cells.Item(76, 8).Value = 2952.3809
$ws.Cells.Item(76, 9).Value = 2894.7368
$ws.Cells.Item(76, 10).Value = 3500
$ws.Cells.Item(76, 11).Value = 2894.7368
$ws.Cells.Item(76, 12).Value = 3500
$ws.Cells.Item(76, 13).Value = -2579.7368
$ws.Cells.Item(76, 14).Value = -4130

$ws.Cells.Item(77, 8).Value = 18430.285
$ws.Cells.Item(77, 9).Value = 18430.285
$ws.Cells.Item(77, 11).Value = 92151.425
$ws.Cells.Item(77, 13).Value = -87471.425

$ws.Cells.Item(79, 8).Value = 2952.3809
$ws.Cells.Item(79, 9).Value = 2894.7368
$ws.Cells.Item(79, 10).Value = 3500
$ws.Cells.Item(79, 11).Value = 2894.7368
$ws.Cells.Item(79, 12).Value = 3500
$ws.Cells.Item(79, 13).Value = -1802.7368
$ws.Cells.Item(79, 14).Value = -5684

$ws.Cells.Item(86, 8).Value = 2791.25
$ws.Cells.Item(86, 10).Value = 2346.2
$ws.Cells.Item(86, 12).Value = 2346.2
$ws.Cells.Item(86, 14).Value = -4592.2

$ws.Cells.Item(89, 8).Value = 2791.25
$ws.Cells.Item(89, 10).Value = 2346.2
$ws.Cells.Item(89, 12).Value = 11731
$ws.Cells.Item(89, 14).Value = -22963

$ws.Cells.Item(98, 8).Value = 6090.9653
$ws.Cells.Item(98, 9).Value = 7726.773
$ws.Cells.Item(98, 11).Value = 7726.773
$ws.Cells.Item(98, 13).Value = -6228.773

$ws.Cells.Item(107, 8).Value = 542.36365
$ws.Cells.Item(107, 9).Value = 500.125
$ws.Cells.Item(107, 11).Value = 500.125
$ws.Cells.Item(107, 13).Value = 1419.875

$ws.Cells.Item(116, 8).Value = 5349
$ws.Cells.Item(116, 10).Value = 5829.2
$ws.Cells.Item(116, 12).Value = 5829.2
$ws.Cells.Item(116, 14).Value = -12713.2

$ws.Cells.Item(122, 8).Value = 6090.9653
$ws.Cells.Item(122, 9).Value = 7726.773
$ws.Cells.Item(122, 11).Value = 23180.319
$ws.Cells.Item(122, 13).Value = -20730.319

$ws.Cells.Item(132, 8).Value = 1873
$ws.Cells.Item(132, 9).Value = 1485
$ws.Cells.Item(132, 10).Value = 3468.111
$ws.Cells.Item(132, 11).Value = 4455
$ws.Cells.Item(132, 12).Value = 10404.333
$ws.Cells.Item(132, 13).Value = -1925
$ws.Cells.Item(132, 14).Value = -15464.333

$ws.Cells.Item(137, 8).Value = 2185.7273
$ws.Cells.Item(137, 9).Value = 2138.2222
$ws.Cells.Item(137, 10).Value = 2399.5
$ws.Cells.Item(137, 11).Value = 6414.6666
$ws.Cells.Item(137, 12).Value = 7198.5
$ws.Cells.Item(137, 13).Value = -3864.6666
$ws.Cells.Item(137, 14).Value = -12298.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 1556.8572
$ws.Cells.Item(26, 9).Value = 979.6
$ws.Cells.Item(26, 11).Value = 979.6
$ws.Cells.Item(26, 13).Value = -649.6

$ws.Cells.Item(32, 8).Value = 2958.6956
$ws.Cells.Item(32, 9).Value = 2913.3333
$ws.Cells.Item(32, 11).Value = 2913.3333
$ws.Cells.Item(32, 13).Value = -2626.3333

$ws.Cells.Item(132, 8).Value = 2170.0588
$ws.Cells.Item(132, 9).Value = 2117.3125
$ws.Cells.Item(132, 11).Value = 6351.9375
$ws.Cells.Item(132, 13).Value = -3821.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 705801.75
$ws.Cells.Item(86, 9).Value = 2807.6667
$ws.Cells.Item(86, 11).Value = 2807.6667
$ws.Cells.Item(86, 13).Value = -1684.6667

$ws.Cells.Item(89, 8).Value = 705801.75
$ws.Cells.Item(89, 9).Value = 2807.6667
$ws.Cells.Item(89, 11).Value = 14038.3335
$ws.Cells.Item(89, 13).Value = -8422.333500000001

$ws.Cells.Item(94, 8).Value = 4775.5713
$ws.Cells.Item(94, 9).Value = 4891.1816
$ws.Cells.Item(94, 10).Value = 4351.6665
$ws.Cells.Item(94, 11).Value = 4891.1816
$ws.Cells.Item(94, 12).Value = 4351.6665
$ws.Cells.Item(94, 13).Value = -4440.1816
$ws.Cells.Item(94, 14).Value = -5253.6665

$ws.Cells.Item(105, 8).Value = 2819.862
$ws.Cells.Item(105, 9).Value = 2438.85
$ws.Cells.Item(105, 10).Value = 3666.5557
$ws.Cells.Item(105, 11).Value = 2438.85
$ws.Cells.Item(105, 12).Value = 3666.5557
$ws.Cells.Item(105, 13).Value = -691.8499999999999
$ws.Cells.Item(105, 14).Value = -7160.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()

$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 13).ClearContents()

$ws.Cells.Item(31, 8).Value = 5284.2
$ws.Cells.Item(31, 9).Value = 2192.1428
$ws.Cells.Item(31, 10).Value = 12499
$ws.Cells.Item(31, 11).Value = 2192.1428
$ws.Cells.Item(31, 12).Value = 12499
$ws.Cells.Item(31, 13).Value = -1897.1428
$ws.Cells.Item(31, 14).Value = -13089

$ws.Cells.Item(34, 8).Value = 5284.2
$ws.Cells.Item(34, 9).Value = 2192.1428
$ws.Cells.Item(34, 10).Value = 12499
$ws.Cells.Item(34, 11).Value = 2192.1428
$ws.Cells.Item(34, 12).Value = 12499
$ws.Cells.Item(34, 13).Value = -1990.1428
$ws.Cells.Item(34, 14).Value = -12903

$ws.Cells.Item(134, 8).Value = 2906.3794
$ws.Cells.Item(134, 9).Value = 2838.1
$ws.Cells.Item(134, 10).Value = 3058.111
$ws.Cells.Item(134, 11).Value = 8514.299999999999
$ws.Cells.Item(134, 12).Value = 9174.332999999999
$ws.Cells.Item(134, 13).Value = -5979.299999999999
$ws.Cells.Item(134, 14).Value = -14244.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 5966.3335
$ws.Cells.Item(112, 10).Value = 1449.5
$ws.Cells.Item(112, 12).Value = 4348.5
$ws.Cells.Item(112, 14).Value = -6564.5

$ws.Cells.Item(132, 8).Value = 2079.739
$ws.Cells.Item(132, 9).Value = 1717.0714
$ws.Cells.Item(132, 10).Value = 2643.889
$ws.Cells.Item(132, 11).Value = 15453.6426
$ws.Cells.Item(132, 12).Value = 23795.001
$ws.Cells.Item(132, 13).Value = -12923.6426
$ws.Cells.Item(132, 14).Value = -28855.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 47624610
$ws.Cells.Item(70, 10).Value = 7249.75
$ws.Cells.Item(70, 12).Value = 7249.75
$ws.Cells.Item(70, 14).Value = -7789.75

$ws.Cells.Item(73, 8).Value = 47624610
$ws.Cells.Item(73, 10).Value = 7249.75
$ws.Cells.Item(73, 12).Value = 7249.75
$ws.Cells.Item(73, 14).Value = -9121.75

$ws.Cells.Item(80, 8).Value = 5915.35
$ws.Cells.Item(80, 9).Value = 5300.375
$ws.Cells.Item(80, 10).Value = 6325.3335
$ws.Cells.Item(80, 11).Value = 5300.375
$ws.Cells.Item(80, 12).Value = 6325.3335
$ws.Cells.Item(80, 13).Value = -4302.375
$ws.Cells.Item(80, 14).Value = -8321.333500000001

$ws.Cells.Item(83, 8).Value = 5915.35
$ws.Cells.Item(83, 9).Value = 5300.375
$ws.Cells.Item(83, 10).Value = 6325.3335
$ws.Cells.Item(83, 11).Value = 26501.875
$ws.Cells.Item(83, 12).Value = 31626.6675
$ws.Cells.Item(83, 13).Value = -21509.875
$ws.Cells.Item(83, 14).Value = -41610.6675

$ws.Cells.Item(122, 8).Value = 12711.192
$ws.Cells.Item(122, 9).Value = 14279.35
$ws.Cells.Item(122, 11).Value = 42838.05
$ws.Cells.Item(122, 13).Value = -40388.05

$ws.Cells.Item(132, 8).Value = 2121.95
$ws.Cells.Item(132, 9).Value = 1880.3704
$ws.Cells.Item(132, 10).Value = 2623.6924
$ws.Cells.Item(132, 11).Value = 5641.1112
$ws.Cells.Item(132, 12).Value = 7871.0772
$ws.Cells.Item(132, 13).Value = -3111.1112
$ws.Cells.Item(132, 14).Value = -12931.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3999
$ws.Cells.Item(46, 10).Value = 3999
$ws.Cells.Item(46, 12).Value = 3999
$ws.Cells.Item(46, 14).Value = -4375

$ws.Cells.Item(100, 8).Value = 7666.5
$ws.Cells.Item(100, 9).Value = 4499.5
$ws.Cells.Item(100, 10).Value = 9250
$ws.Cells.Item(100, 11).Value = 4499.5
$ws.Cells.Item(100, 12).Value = 9250
$ws.Cells.Item(100, 13).Value = -3958.5
$ws.Cells.Item(100, 14).Value = -10332

$ws.Cells.Item(132, 8).Value = 33336160
$ws.Cells.Item(132, 9).Value = 2873.818
$ws.Cells.Item(132, 10).Value = 125002700
$ws.Cells.Item(132, 11).Value = 8621.454000000002
$ws.Cells.Item(132, 12).Value = 375008100
$ws.Cells.Item(132, 13).Value = -6091.454000000002
$ws.Cells.Item(132, 14).Value = -375013160

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 76664.664
$ws.Cells.Item(27, 10).Value = 76664.664
$ws.Cells.Item(27, 12).Value = 76664.664
$ws.Cells.Item(27, 14).Value = -76802.664

$ws.Cells.Item(81, 8).Value = 4911.1665
$ws.Cells.Item(81, 10).Value = 5590.8
$ws.Cells.Item(81, 12).Value = 11181.6
$ws.Cells.Item(81, 14).Value = -13303.6

$ws.Cells.Item(84, 8).Value = 4911.1665
$ws.Cells.Item(84, 10).Value = 5590.8
$ws.Cells.Item(84, 12).Value = 55908
$ws.Cells.Item(84, 14).Value = -66516

$ws.Cells.Item(126, 8).Value = 1892.1765
$ws.Cells.Item(126, 9).Value = 1804.3636
$ws.Cells.Item(126, 11).Value = 5413.0908
$ws.Cells.Item(126, 13).Value = -2943.0908

$ws.Cells.Item(132, 8).Value = 1965.2778
$ws.Cells.Item(132, 9).Value = 1537.5714
$ws.Cells.Item(132, 10).Value = 3462.25
$ws.Cells.Item(132, 11).Value = 4612.7142
$ws.Cells.Item(132, 12).Value = 10386.75
$ws.Cells.Item(132, 13).Value = -2082.7142
$ws.Cells.Item(132, 14).Value = -15446.75
